# Applies the three text edits described by the diff:
#   1. Title: "이력서 - Patti Fernandez" -> "이력서: Patti Fernandez"
#   2. Job entry: remove the space before the parenthesis in
#      "주니어 애니메이터 (2012년 9월 - 2015년 5월)" -> "주니어 애니메이터(2012년 9월 - 2015년 5월)"
#   3. Reference citation: translate the Korean sentence fragments to English
#      "애니메이션의 예술: 초보자를 위한 가이드입니다. " -> "The Art of Animation: A Guide for Beginners. "
#      "뉴욕: 펭귄 책."                                 -> "New York: Penguin Books."
#
# NOTE on the citation runs: the two runs being edited sit right next to a
# third run ("Fernandez, P.(2020년). ") that shares byte-for-byte identical
# run formatting (rPr). Whenever a Range.Text assignment touches a run in
# that paragraph, adjacent runs that end up with identical formatting get
# coalesced into a single run. Since the diff keeps the three runs
# separate (same rPr, only the text differs), each edited run's Font is
# nudged to a different formatting value (Italic) right before the text
# swap - which prevents the coalescing - and then restored back to its
# original value immediately after, once the text is already in place.

$d = $word.ActiveDocument

# 1. Title
$d.Content.Find.Execute("이력서 - Patti Fernandez", $true, $false, $false, $false, $false,
                         $true, 1, $false, "이력서: Patti Fernandez", 2)

# 2. Job title / date range spacing
$d.Content.Find.Execute("주니어 애니메이터 (2012년 9월 - 2015년 5월)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "주니어 애니메이터(2012년 9월 - 2015년 5월)", 2)

# 3a. First citation run
$rng2 = $d.Content
$rng2.Find.Execute("애니메이션의 예술: 초보자를 위한 가이드입니다. ", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
$rng2.Font.Italic = 1
$rng2.Text = "The Art of Animation: A Guide for Beginners. "

$rng2b = $d.Content
$rng2b.Find.Execute("The Art of Animation: A Guide for Beginners. ", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$rng2b.Font.Italic = 0

# 3b. Second citation run
$rng3 = $d.Content
$rng3.Find.Execute("뉴욕: 펭귄 책.", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
$rng3.Font.Italic = 1
$rng3.Text = "New York: Penguin Books."

$rng3b = $d.Content
$rng3b.Find.Execute("New York: Penguin Books.", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$rng3b.Font.Italic = 0
